# Updated cryptos list values (Price column D, Volume(1h) column E,
# plus a Coin/Link/Price/Volume swap between rows 48 and 49) per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $text) {
    # Force the cell to stay plain text like the original inline-string cell,
    # instead of Excel auto-coercing a numeric-looking string (e.g. "211.06")
    # into a real number; then restore the default style so no stray
    # number-format/style index is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}
Set-TextCell $ws.Range('D2') '26.095.41'
$ws.Range('E2').Value = '  -1.08%  '
Set-TextCell $ws.Range('D3') '1.676.92'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  -0.46%  '
Set-TextCell $ws.Range('D5') '210.98'
$ws.Range('E5').Value = '  -3.46%  '
Set-TextCell $ws.Range('D6') '0.5273'
$ws.Range('E6').Value = '  -4.68%  '
$ws.Range('E7').Value = '  -0.40%  '
Set-TextCell $ws.Range('D8') '0.2678'
$ws.Range('E8').Value = '  -1.21%  '
Set-TextCell $ws.Range('D9') '0.06311'
$ws.Range('E9').Value = '  -2.97%  '
Set-TextCell $ws.Range('D10') '21.27'
$ws.Range('E10').Value = '  -3.89%  '
Set-TextCell $ws.Range('D11') '0.07569'
$ws.Range('E11').Value = '  -0.11%  '
Set-TextCell $ws.Range('D12') '1.675.76'
$ws.Range('E12').Value = '  -0.65%  '
Set-TextCell $ws.Range('D13') '4.508'
$ws.Range('E13').Value = '  -1.03%  '
Set-TextCell $ws.Range('D14') '0.5684'
$ws.Range('E14').Value = '  -2.32%  '
Set-TextCell $ws.Range('D15') '0.000008135'
$ws.Range('E15').Value = '  -3.96%  '
Set-TextCell $ws.Range('D16') '65.71'
$ws.Range('E16').Value = '  +0.38%  '
Set-TextCell $ws.Range('D17') '26.127.97'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('E18').Value = '  -0.37%  '
Set-TextCell $ws.Range('D19') '4.857'
$ws.Range('E19').Value = '  -1.75%  '
Set-TextCell $ws.Range('D20') '10.61'
$ws.Range('E20').Value = '  -2.92%  '
Set-TextCell $ws.Range('D21') '189.37'
$ws.Range('E21').Value = '  -1.05%  '
Set-TextCell $ws.Range('D22') '6.196'
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('E23').Value = '  -0.40%  '
Set-TextCell $ws.Range('D24') '148.54'
$ws.Range('E24').Value = '  -0.14%  '
Set-TextCell $ws.Range('D25') '0.1254'
$ws.Range('E25').Value = '  -5.37%  '
Set-TextCell $ws.Range('D26') '7.651'
$ws.Range('E26').Value = '  -3.24%  '
Set-TextCell $ws.Range('D27') '16.13'
$ws.Range('E27').Value = '  +1.86%  '
Set-TextCell $ws.Range('D28') '0.06352'
$ws.Range('E28').Value = '  +0.45%  '
Set-TextCell $ws.Range('D29') '1.354'
$ws.Range('E29').Value = '  -3.10%  '
Set-TextCell $ws.Range('D30') '1.285'
$ws.Range('E30').Value = '  -3.11%  '
Set-TextCell $ws.Range('D31') '3.539'
$ws.Range('E31').Value = '  -1.65%  '
Set-TextCell $ws.Range('D32') '3.537'
$ws.Range('E32').Value = '  -1.30%  '
Set-TextCell $ws.Range('D33') '1.671'
$ws.Range('E33').Value = '  -0.15%  '
Set-TextCell $ws.Range('D34') '1.011'
$ws.Range('E34').Value = '  -2.97%  '
Set-TextCell $ws.Range('D35') '0.6058'
$ws.Range('E35').Value = '  -2.81%  '
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('E37').Value = '  -0.12%  '
Set-TextCell $ws.Range('D38') '6.151'
$ws.Range('E38').Value = '  -1.39%  '
$ws.Range('E39').Value = '  -1.21%  '
Set-TextCell $ws.Range('D40') '1.099.68'
$ws.Range('E40').Value = '  -1.55%  '
Set-TextCell $ws.Range('D41') '0.8725'
$ws.Range('E41').Value = '  -0.56%  '
Set-TextCell $ws.Range('D42') '1.005'
$ws.Range('E42').Value = '  -0.94%  '
$ws.Range('E43').Value = '  -0.51%  '
Set-TextCell $ws.Range('D44') '1.828.17'
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('E45').Value = '  +0.11%  '
Set-TextCell $ws.Range('D46') '57.08'
$ws.Range('E46').Value = '  -0.58%  '
Set-TextCell $ws.Range('D47') '1.006'
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws.Range('D48') '0.05252'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws.Range('D49') '7.989'
$ws.Range('E49').Value = '  -2.77%  '
$ws.Range('E50').Value = '  -0.77%  '
Set-TextCell $ws.Range('D51') '5.954'
$ws.Range('E51').Value = '  -2.19%  '
